$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1), styled like the other headers (copy style from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 7
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7
